$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the NA US (test4ca -> test18) environment values in row 2.
$ws.Range("A2").Value = "https://test18.cliotest.com/backoffice/control/main"
$ws.Range("C2").Value = "https://test18.cliotest.com/cabicentral/control/main"
$ws.Range("D2").Value = "https://test18.cliotest.com/warehouse/control/main"
$ws.Range("F2").Value = "virtual_cabitest18"
$ws.Range("G2").Value = "test18"
